# Applies the MAINE_2019.xlsx data-cleaning fixes:
#  - rename header columns to the cleaned snake_case names
#  - normalize "de la"/"de" -> "De La"/"De" capitalization in a few place names
#  - tiny floating point precision correction on a few percentage cells
#  - drop the trailing footnote/metadata rows (33-37), shrinking the used range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Text capitalization fixes ------------------------------------------
$ws.Range("B2").Value  = "Amatenango De La Frontera"
$ws.Range("B5").Value  = "Mazapa De Madero"
$ws.Range("A9").Value  = "Estado De México"
$ws.Range("B16").Value = "San Cristóbal De La Barranca"

# --- Floating point precision corrections --------------------------------
$ws.Range("D13").Value = 0.09677419354838708
$ws.Range("D14").Value = 0.09677419354838708
$ws.Range("D17").Value = 0.09677419354838708

# --- Remove trailing metadata rows (33-37) --------------------------------
$ws.Range("A33:A37").EntireRow.Delete()
